$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the table (ListObject) first to include the new column (I)
$table = $ws.ListObjects.Item("Tablo1")
$table.Resize($ws.Range("A1:I29"))

# Set the new column header text (this also syncs the table's column name)
$ws.Cells.Item(1, 9).Value = "Sınıf Sayısı"

# Fill column I (rows 2-29) with value 1
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# Update selection to match target state
$ws.Range("Q15").Select()
